$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (time changed from 04:05 to 05:22)
$ws.Range("A1").Value = "Datos actualizados a 2 de Septiembre de 2020 a las 05:22"

# Row 31 - Bolivia
$ws.Range("B31").Value = 117267
$ws.Range("C31").Value = 669
$ws.Range("D31").Value = 62124
$ws.Range("E31").Value = 50042
$ws.Range("G31").Value = 74
$ws.Range("H31").Value = 5101

# Row 33 - Kazajistan
$ws.Range("B33").Value = 105944
$ws.Range("C33").Value = 72
$ws.Range("D33").Value = 97371
$ws.Range("E33").Value = 6985

# Row 40 - Belgica
$ws.Range("B40").Value = 85487
$ws.Range("C40").Value = 251
$ws.Range("D40").Value = 18457
$ws.Range("E40").Value = 57133
$ws.Range("G40").Value = 2
$ws.Range("H40").Value = 9897

# Row 50 - Honduras
$ws.Range("B50").Value = 61769
$ws.Range("C50").Value = 755
$ws.Range("D50").Value = 10430
$ws.Range("E50").Value = 49451
$ws.Range("G50").Value = 15
$ws.Range("H50").Value = 1888

# Row 175 - San Martin (Parte Holandesa)
$ws.Range("B175").Value = 482
$ws.Range("C175").Value = 6
$ws.Range("D175").Value = 241
$ws.Range("E175").Value = 222
$ws.Range("G175").Value = 2
$ws.Range("H175").Value = 19
